# Updated cryptos list on Mon Jun 10 14:53:17 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# coin rows, and swaps the FirstDigitalUSD / Monero rows (40 <-> 41) to
# reflect the new ranking order.
#
# Every value in this sheet is stored as text (inline strings), including
# ones that look numeric ("1.00", "0.650", ...). A plain
# `$ws.Range(...).Value = "1.00"` assignment lets Excel's type inference
# kick in and silently reinterpret the text as a number, which strips
# meaningful trailing zeros (e.g. "1.00" -> "1", "0.650" -> "0.65"). To
# keep the exact original text representation, we briefly force the
# cell's number format to Text ("@") before writing the value, then
# restore the cell's original style so no visible formatting changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" "69.605.05"
Set-TextValue "E2" "  +0.19%  "
Set-TextValue "D3" "3.677.67"
Set-TextValue "E3" "  -0.35%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.08%  "
Set-TextValue "D5" "647.48"
Set-TextValue "E5" "  -4.72%  "
Set-TextValue "D6" "159.83"
Set-TextValue "E6" "  -0.93%  "
Set-TextValue "D8" "0.501"
Set-TextValue "E8" "  +1.15%  "
Set-TextValue "E9" "  -0.97%  "
Set-TextValue "D10" "7.19"
Set-TextValue "E10" "  +0.22%  "
Set-TextValue "D11" "0.443"
Set-TextValue "E11" "  +0.67%  "
Set-TextValue "D12" "0.0000232"
Set-TextValue "E12" "  -0.84%  "
Set-TextValue "D13" "4.297.70"
Set-TextValue "E13" "  -0.33%  "
Set-TextValue "D14" "32.72"
Set-TextValue "E14" "  +0.66%  "
Set-TextValue "D15" "3.675.45"
Set-TextValue "E15" "  -0.50%  "
Set-TextValue "D16" "69.541.47"
Set-TextValue "E16" "  +0.18%  "
Set-TextValue "E17" "  +0.15%  "
Set-TextValue "D18" "16.03"
Set-TextValue "E18" "  -0.14%  "
Set-TextValue "D19" "6.50"
Set-TextValue "E19" "  +0.14%  "
Set-TextValue "D20" "469.42"
Set-TextValue "E20" "  -0.51%  "
Set-TextValue "D21" "10.07"
Set-TextValue "E21" "  +2.74%  "
Set-TextValue "D22" "0.650"
Set-TextValue "E22" "  -0.23%  "
Set-TextValue "D23" "79.65"
Set-TextValue "E23" "  -0.94%  "
Set-TextValue "D24" "3.820.46"
Set-TextValue "E24" "  -0.41%  "
Set-TextValue "E25" "  +0.03%  "
Set-TextValue "D26" "0.0000125"
Set-TextValue "E26" "  -0.92%  "
Set-TextValue "D27" "10.99"
Set-TextValue "E27" "  +0.86%  "
Set-TextValue "D28" "9.11"
Set-TextValue "E28" "  -0.47%  "
Set-TextValue "D29" "2.64"
Set-TextValue "E29" "  -2.49%  "
Set-TextValue "D30" "1.72"
Set-TextValue "E30" "  -1.67%  "
Set-TextValue "D31" "2.01"
Set-TextValue "E31" "  -0.18%  "
Set-TextValue "D32" "1.00"
Set-TextValue "E32" "  +0.29%  "
Set-TextValue "D33" "26.76"
Set-TextValue "E33" "  -0.91%  "
Set-TextValue "D34" "6.45"
Set-TextValue "E34" "  -2.09%  "
Set-TextValue "D35" "3.671.10"
Set-TextValue "E35" "  -0.24%  "
Set-TextValue "E36" "  +0.00%  "
Set-TextValue "D37" "8.47"
Set-TextValue "E37" "  +0.15%  "
Set-TextValue "D39" "5.89"
Set-TextValue "E39" "  -5.19%  "
Set-TextValue "B40" "FirstDigitalUSD"
Set-TextValue "C40" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D40" "0.999"
Set-TextValue "E40" "  -0.05%  "
Set-TextValue "B41" "Monero"
Set-TextValue "C41" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D41" "176.34"
Set-TextValue "E41" "  +4.62%  "
Set-TextValue "D42" "2.22"
Set-TextValue "E42" "  -1.87%  "
Set-TextValue "D43" "0.0901"
Set-TextValue "E43" "  -0.37%  "
Set-TextValue "D44" "0.927"
Set-TextValue "E44" "  -1.61%  "
Set-TextValue "D45" "47.17"
Set-TextValue "E45" "  +1.05%  "
Set-TextValue "D46" "28.93"
Set-TextValue "E46" "  +2.87%  "
Set-TextValue "D47" "2.70"
Set-TextValue "E47" "  -1.33%  "
Set-TextValue "E48" "  -1.42%  "
Set-TextValue "E49" "  -5.13%  "
Set-TextValue "D50" "7.83"
Set-TextValue "E50" "  -0.91%  "
Set-TextValue "D51" "1.24"
Set-TextValue "E51" "  -3.62%  "
